{"js": "// Remove the first (title/heading) paragraph of the document body:\n// \"Deformable Sensors based on Architectured 2D Materials\" (Heading1 style).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  const firstParagraph = paragraphs.items[0];\n  firstParagraph.load(\"text\");\n  await context.sync();\n\n  if (firstParagraph.text.indexOf(\"Deformable Sensors based on\") !== -1) {\n    firstParagraph.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the first (title/heading) paragraph of the document body:\n# \"Deformable Sensors based on Architectured 2D Materials\" (Heading1 style).\n$d = $word.ActiveDocument\n\n$firstParagraph = $d.Paragraphs.First\nif ($firstParagraph.Range.Text -like \"*Deformable Sensors based on*\") {\n    $firstParagraph.Range.Delete()\n}\n"}
